$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-48 with new values
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 6
$ws.Cells.Item(2, 3).Value = 58
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 7
$ws.Cells.Item(3, 3).Value = 32
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 8
$ws.Cells.Item(4, 3).Value = 27
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 9
$ws.Cells.Item(5, 3).Value = 13
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 10
$ws.Cells.Item(6, 3).Value = 12
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 12
$ws.Cells.Item(7, 3).Value = 11
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 11
$ws.Cells.Item(8, 3).Value = 10
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 13
$ws.Cells.Item(9, 3).Value = 10
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 27
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 14
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 15
$ws.Cells.Item(12, 3).Value = 5
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 31
$ws.Cells.Item(13, 3).Value = 4
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 18
$ws.Cells.Item(14, 3).Value = 4
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 19
$ws.Cells.Item(15, 3).Value = 4
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 30
$ws.Cells.Item(16, 3).Value = 3
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 25
$ws.Cells.Item(17, 3).Value = 3
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 26
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 32
$ws.Cells.Item(19, 3).Value = 3
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = 23
$ws.Cells.Item(20, 3).Value = 3
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 22
$ws.Cells.Item(21, 3).Value = 3
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 20
$ws.Cells.Item(22, 3).Value = 3
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 17
$ws.Cells.Item(23, 3).Value = 3
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = 34
$ws.Cells.Item(24, 3).Value = 2
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = 41
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = 48
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 33
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = 37
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 130
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 113
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = 87
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = 85
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = 84
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = 83
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = 74
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = 70
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = 68
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = 67
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).Value = 58
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).Value = 55
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).Value = 51
$ws.Cells.Item(41, 3).Value = 1
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).Value = 28
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).Value = 46
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = 44
$ws.Cells.Item(44, 3).Value = 1
$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 2).Value = 43
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 42
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 16
$ws.Cells.Item(47, 3).Value = 1
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 264
$ws.Cells.Item(48, 3).Value = 1

# Remove now-unused rows 49-72 (delete bottom-up)
for ($r = 72; $r -ge 49; $r--) {
    $ws.Rows.Item($r).Delete()
}
